$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-54 down to 51-55
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with the new record's data
$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value = "La Araucanía"
$ws.Cells.Item(50, 4).Value = 44776
$ws.Cells.Item(50, 5).Value = 9
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100108
$ws.Cells.Item(50, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value = 100108003
$ws.Cells.Item(50, 10).Value = "Maracuyá"
$ws.Cells.Item(50, 11).Value = "Sin especificar"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 40
$ws.Cells.Item(50, 14).Value = 35000
$ws.Cells.Item(50, 15).Value = 35000
$ws.Cells.Item(50, 16).Value = 35000
$ws.Cells.Item(50, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(50, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 19).Value = 1944
$ws.Cells.Item(50, 20).Value = 18
